$wb = $excel.ActiveWorkbook

# Sheet: 展览 (Exhibitions)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 409
$ws1.Range("F5").Value = 1231
$ws1.Range("F6").Value = 462
$ws1.Range("F7").Value = 7461
$ws1.Range("F10").Value = 2071
$ws1.Range("F11").Value = 8117
$ws1.Range("F14").Value = 5555
$ws1.Range("F16").Value = 2506
$ws1.Range("F17").Value = 1071
$ws1.Range("F18").Value = 4572
$ws1.Range("F22").Value = 21
$ws1.Range("F23").Value = 437
$ws1.Range("F24").Value = 1316
$ws1.Range("F26").Value = 2581
$ws1.Range("F28").Value = 298
$ws1.Range("F30").Value = 219
$ws1.Range("F32").Value = 9
$ws1.Range("F34").Value = 1577
$ws1.Range("F35").Value = 36
$ws1.Range("F37").Value = 2498
$ws1.Range("F38").Value = 2250

# Sheet: 演出 (Performances)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 89
$ws2.Range("F4").Value = 28

# Sheet: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F5").Value = 409
$ws4.Range("F6").Value = 1231
$ws4.Range("F7").Value = 462
$ws4.Range("F8").Value = 7461
$ws4.Range("F11").Value = 2071
$ws4.Range("F12").Value = 8117
$ws4.Range("F15").Value = 5555
$ws4.Range("F17").Value = 2506
$ws4.Range("F18").Value = 1071
$ws4.Range("F19").Value = 4572
$ws4.Range("F24").Value = 21
$ws4.Range("F25").Value = 89
$ws4.Range("F26").Value = 437
$ws4.Range("F27").Value = 1316
$ws4.Range("F29").Value = 2581
$ws4.Range("F31").Value = 298
$ws4.Range("F33").Value = 219
$ws4.Range("F34").Value = 28
$ws4.Range("F35").Value = 616
$ws4.Range("F36").Value = 9
$ws4.Range("F39").Value = 1577
$ws4.Range("F40").Value = 36
$ws4.Range("F42").Value = 2498
$ws4.Range("F44").Value = 2250
